$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1, Col 1: "61÷7=8, 5" -> "55÷7=7, 6"
$cell = $t.Cell(1, 1)
$cell.Range.Text = "55÷7=7, 6"

# Row 1, Col 2: "88÷6=14, 4" -> "55÷6=9, 1"
$cell = $t.Cell(1, 2)
$cell.Range.Text = "55÷6=9, 1"

# Row 1, Col 3: "99÷4=24, 3" -> "56÷9=6, 2"
$cell = $t.Cell(1, 3)
$cell.Range.Text = "56÷9=6, 2"

# Row 1, Col 4: "75÷7=10, 5" -> "10÷6=1, 4"
$cell = $t.Cell(1, 4)
$cell.Range.Text = "10÷6=1, 4"

# Row 1, Col 5: "41÷8=5, 1" -> "11÷5=2, 1"
$cell = $t.Cell(1, 5)
$cell.Range.Text = "11÷5=2, 1"

# Row 5, Col 1: "74÷7=10, 4" -> "54÷8=6, 6"
$cell = $t.Cell(5, 1)
$cell.Range.Text = "54÷8=6, 6"

# Row 5, Col 2: "67÷6=11, 1" -> "49÷6=8, 1"
$cell = $t.Cell(5, 2)
$cell.Range.Text = "49÷6=8, 1"

# Row 5, Col 3: "49÷4=12, 1" -> "66÷2=33, 0"
$cell = $t.Cell(5, 3)
$cell.Range.Text = "66÷2=33, 0"

# Row 5, Col 4: "55÷5=11, 0" -> "29÷4=7, 1"
$cell = $t.Cell(5, 4)
$cell.Range.Text = "29÷4=7, 1"

# Row 5, Col 5: "39÷9=4, 3" -> "76÷4=19, 0"
$cell = $t.Cell(5, 5)
$cell.Range.Text = "76÷4=19, 0"

# Row 9, Col 1: "52÷5=10, 2" -> "98÷4=24, 2"
$cell = $t.Cell(9, 1)
$cell.Range.Text = "98÷4=24, 2"

# Row 9, Col 2: "10÷6=1, 4" -> "14÷7=2, 0"
$cell = $t.Cell(9, 2)
$cell.Range.Text = "14÷7=2, 0"

# Row 9, Col 3: "37÷4=9, 1" -> "32÷4=8, 0"
$cell = $t.Cell(9, 3)
$cell.Range.Text = "32÷4=8, 0"

# Row 9, Col 4: "79÷8=9, 7" -> "15÷6=2, 3"
$cell = $t.Cell(9, 4)
$cell.Range.Text = "15÷6=2, 3"

# Row 9, Col 5: "54÷9=6, 0" -> "48÷6=8, 0"
$cell = $t.Cell(9, 5)
$cell.Range.Text = "48÷6=8, 0"

# Row 13, Col 1: "21÷4=5, 1" -> "23÷2=11, 1"
$cell = $t.Cell(13, 1)
$cell.Range.Text = "23÷2=11, 1"

# Row 13, Col 2: "33÷4=8, 1" -> "55÷4=13, 3"
$cell = $t.Cell(13, 2)
$cell.Range.Text = "55÷4=13, 3"

# Row 13, Col 3: "13÷4=3, 1" -> "85÷2=42, 1"
$cell = $t.Cell(13, 3)
$cell.Range.Text = "85÷2=42, 1"

# Row 13, Col 4: "87÷2=43, 1" -> "63÷6=10, 3"
$cell = $t.Cell(13, 4)
$cell.Range.Text = "63÷6=10, 3"

# Row 13, Col 5: "40÷3=13, 1" -> "24÷3=8, 0"
$cell = $t.Cell(13, 5)
$cell.Range.Text = "24÷3=8, 0"

# Row 17, Col 1: "72÷2=36, 0" -> "91÷7=13, 0"
$cell = $t.Cell(17, 1)
$cell.Range.Text = "91÷7=13, 0"

# Row 17, Col 2: "59÷2=29, 1" -> "41÷6=6, 5"
$cell = $t.Cell(17, 2)
$cell.Range.Text = "41÷6=6, 5"

# Row 17, Col 3: "92÷8=11, 4" -> "62÷7=8, 6"
$cell = $t.Cell(17, 3)
$cell.Range.Text = "62÷7=8, 6"

# Row 17, Col 4: "35÷6=5, 5" -> "13÷7=1, 6"
$cell = $t.Cell(17, 4)
$cell.Range.Text = "13÷7=1, 6"

# Row 17, Col 5: "58÷2=29, 0" -> "50÷6=8, 2"
$cell = $t.Cell(17, 5)
$cell.Range.Text = "50÷6=8, 2"
